# Update cosinor-per-day analysis results (CircaDB / CircadiPy simulation re-run).
# Values below were recomputed for sawtooth_0.1 fixed-period-9 cosinor analysis.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.7109059811527194
$ws.Range("I2").Value = 0.7109059811527194
$ws.Range("L2").Value = 3.155673832846288
$ws.Range("M2").Value = "[-7.163826612312384, 13.47517427800496]"
$ws.Range("N2").Value = 0.5410588463476371
$ws.Range("O2").Value = 0.5410588463476371
$ws.Range("P2").Value = -0.69184222601577
$ws.Range("Q2").Value = "[-3.8240006674326206, 2.4403162154010807]"
$ws.Range("R2").Value = 0.6585374292306385
$ws.Range("S2").Value = 0.6585374292306385
$ws.Range("T2").Value = 14.43448211463249
$ws.Range("U2").Value = "[8.905947884412932, 19.963016344852054]"
$ws.Range("V2").Value = 0.000003870457323307264
$ws.Range("W2").Value = 0.000003870457323307264
$ws.Range("X2").Value = 2.86176176176183
$ws.Range("Y2").Value = -10.09421421421446
$ws.Range("Z2").Value = 15.81773773773812
# Row 3
$ws.Range("B3").Value = 1
$ws.Range("H3").Value = 0.04076396158379902
$ws.Range("I3").Value = 0.04076396158379902
$ws.Range("L3").Value = 7.623535835370946
$ws.Range("M3").Value = "[0.0191229782517901, 15.227948692490102]"
$ws.Range("N3").Value = 0.04945062882090778
$ws.Range("O3").Value = 0.04945062882090778
$ws.Range("P3").Value = -1.94973718240808
$ws.Range("Q3").Value = "[-3.4340532309510046, -0.4654211338651546]"
$ws.Range("R3").Value = 0.01119148083301558
$ws.Range("S3").Value = 0.01119148083301558
$ws.Range("T3").Value = 12.85329056456845
$ws.Range("U3").Value = "[8.660955707643197, 17.045625421493696]"
$ws.Range("V3").Value = 0.000000171568766749175
$ws.Range("W3").Value = 0.000000171568766749175
$ws.Range("X3").Value = 8.06496496496516
$ws.Range("Y3").Value = 1.925185185185234
$ws.Range("Z3").Value = 14.20474474474508
# Row 4
$ws.Range("H4").Value = 0.7383299622524367
$ws.Range("I4").Value = 0.7383299622524367
$ws.Range("L4").Value = 3.111591290551689
$ws.Range("M4").Value = "[-7.690282177055079, 13.913464758158458]"
$ws.Range("N4").Value = 0.5646850458635027
$ws.Range("O4").Value = 0.5646850458635027
$ws.Range("P4").Value = -1.685579241565695
$ws.Range("Q4").Value = "[-4.824027157764508, 1.4528686746331179]"
$ws.Range("R4").Value = 0.2851385735604794
$ws.Range("S4").Value = 0.2851385735604794
$ws.Range("T4").Value = 17.06270539881688
$ws.Range("U4").Value = "[11.208426087727293, 22.916984709906476]"
$ws.Range("V4").Value = 0.0000004863517841791065
$ws.Range("W4").Value = 0.0000004863517841791065
$ws.Range("X4").Value = 6.972292292292462
$ws.Range("Y4").Value = -6.009699699699844
$ws.Range("Z4").Value = 19.95428428428477
# Row 5
$ws.Range("H5").Value = 0.7590539259104988
$ws.Range("I5").Value = 0.7590539259104988
$ws.Range("L5").Value = 2.63010279743795
$ws.Range("M5").Value = "[-6.997031379040538, 12.257236973916438]"
$ws.Range("N5").Value = 0.5848713205197955
$ws.Range("O5").Value = 0.5848713205197955
$ws.Range("P5").Value = -2.012631930227695
$ws.Range("Q5").Value = "[-5.151079846426507, 1.1258159859711174]"
$ws.Range("R5").Value = 0.203087331125378
$ws.Range("S5").Value = 0.203087331125378
$ws.Range("T5").Value = 14.93813105146852
$ws.Range("U5").Value = "[9.869570997867825, 20.006691105069216]"
$ws.Range("V5").Value = 0.0000003885656776336788
$ws.Range("W5").Value = 0.0000003885656776336788
$ws.Range("X5").Value = 8.325125125125322
$ws.Range("Y5").Value = -4.656866866866981
$ws.Range("Z5").Value = 21.30711711711762
# Row 6
$ws.Range("F6").Value = 22.86000000000013
$ws.Range("H6").Value = 0.3567271136453359
$ws.Range("I6").Value = 0.3567271136453359
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = 4.159728135172344
$ws.Range("M6").Value = "[-3.0570956535816283, 11.376551923926316]"
$ws.Range("N6").Value = 0.2517977923193542
$ws.Range("O6").Value = 0.2517977923193542
$ws.Range("P6").Value = -3.069263693597235
$ws.Range("Q6").Value = "[-6.207711609796047, 0.06918422260157708]"
$ws.Range("R6").Value = 0.05504622632488809
$ws.Range("S6").Value = 0.05504622632488809
$ws.Range("T6").Value = 11.52758037854033
$ws.Range("U6").Value = "[7.363332621476911, 15.691828135603743]"
$ws.Range("V6").Value = 0.000001326437543580639
$ws.Range("W6").Value = 0.000001326437543580639
$ws.Range("X6").Value = 11.16684684684691
$ws.Range("Y6").Value = -0.2517117117117138
$ws.Range("Z6").Value = 22.58540540540553
# Row 7
$ws.Range("B7").Value = 1
$ws.Range("F7").Value = 22.86000000000013
$ws.Range("H7").Value = 0.01092753114853739
$ws.Range("I7").Value = 0.01092753114853739
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = 10.34053816623389
$ws.Range("M7").Value = "[2.212026248775807, 18.469050083691968]"
$ws.Range("N7").Value = 0.01381913970804138
$ws.Range("O7").Value = 0.01381913970804138
$ws.Range("P7").Value = -2.956053147521927
$ws.Range("Q7").Value = "[-4.075579658711083, -1.836526636332772]"
$ws.Range("R7").Value = 0.000003167789150992917
$ws.Range("S7").Value = 0.000003167789150992917
$ws.Range("T7").Value = 16.01325918378441
$ws.Range("U7").Value = "[11.266227627456743, 20.76029074011207]"
$ws.Range("V7").Value = 0.00000002058845627317396
$ws.Range("W7").Value = 0.00000002058845627317396
$ws.Range("X7").Value = 10.75495495495502
$ws.Range("Y7").Value = 6.681801801801839
$ws.Range("Z7").Value = 14.82810810810819
# Row 8
$ws.Range("F8").Value = 22.86000000000013
$ws.Range("H8").Value = 0.7218370469745605
$ws.Range("I8").Value = 0.7218370469745605
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = 2.253837207735103
$ws.Range("M8").Value = "[-5.097330489410172, 9.605004904880378]"
$ws.Range("N8").Value = 0.540007122982259
$ws.Range("O8").Value = 0.540007122982259
$ws.Range("P8").Value = 3.037816319687427
$ws.Range("Q8").Value = "[-0.08805264694746207, 6.163685286322317]"
$ws.Range("R8").Value = 0.05652397561660916
$ws.Range("S8").Value = 0.05652397561660916
$ws.Range("T8").Value = 10.19001389453205
$ws.Range("U8").Value = "[6.142644445606528, 14.237383343457577]"
$ws.Range("V8").Value = 0.000007262445536060369
$ws.Range("W8").Value = 0.000007262445536060369
$ws.Range("X8").Value = 11.80756756756764
$ws.Range("Y8").Value = 0.4347747747747768
$ws.Range("Z8").Value = 23.1803603603605
# Row 9
$ws.Range("F9").Value = 22.86000000000013
$ws.Range("H9").Value = 0.4782253801701852
$ws.Range("I9").Value = 0.4782253801701852
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = 3.522974311910819
$ws.Range("M9").Value = "[-4.02603437565126, 11.0719829994729]"
$ws.Range("N9").Value = 0.3522653563627116
$ws.Range("O9").Value = 0.3522653563627116
$ws.Range("P9").Value = -2.012631930227695
$ws.Range("Q9").Value = "[-5.151079846426507, 1.1258159859711174]"
$ws.Range("R9").Value = 0.203087331125378
$ws.Range("S9").Value = 0.203087331125378
$ws.Range("T9").Value = 12.43333192504883
$ws.Range("U9").Value = "[8.402465725757096, 16.464198124340555]"
$ws.Range("V9").Value = 0.0000001508884979184444
$ws.Range("W9").Value = 0.0000001508884979184444
$ws.Range("X9").Value = 7.322522522522563
$ws.Range("Y9").Value = -4.096036036036061
$ws.Range("Z9").Value = 18.74108108108119
# Row 10
$ws.Range("F10").Value = 22.86000000000013
$ws.Range("H10").Value = 0.8374859948050791
$ws.Range("I10").Value = 0.8374859948050791
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = 2.287364830771431
$ws.Range("M10").Value = "[-7.669750028773936, 12.244479690316798]"
$ws.Range("N10").Value = 0.645821145788406
$ws.Range("O10").Value = 0.645821145788406
$ws.Range("P10").Value = -1.647842392873925
$ws.Range("Q10").Value = "[-4.786290309072737, 1.490605523324887]"
$ws.Range("R10").Value = 0.2959251409707098
$ws.Range("S10").Value = 0.2959251409707098
$ws.Range("T10").Value = 14.45751623646423
$ws.Range("U10").Value = "[9.11962184875211, 19.795410624176355]"
$ws.Range("V10").Value = 0.000001994488884271561
$ws.Range("W10").Value = 0.000001994488884271561
$ws.Range("X10").Value = 5.995315315315349
$ws.Range("Y10").Value = -5.423243243243276
$ws.Range("Z10").Value = 17.41387387387397
